# This script appends newly-collected breakout alert rows to the three data
# sheets of the workbook (three_line, two_line, ph_pl_breakout_line), matching
# the 'break out stock.yaml completed' run for 06-06-2024. Date-valued cells are
# written as Excel serial numbers with the same custom date/time NumberFormat
# already used elsewhere in each sheet.

$wb = $excel.ActiveWorkbook
$DATE_FMT = "YYYY-MM-DD HH:MM:SS"

# --- three_line: append rows 66-69 ---
$ws = $wb.Worksheets("three_line")
$newRows = @(
    @{Row=66; Cells=@(@{Col='A'; Value="PRSMJOHNSN.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45406.38541666666; IsDate=$true}; @{Col='E'; Value=171.3500061035156; IsDate=$false}; @{Col='F'; Value=45412.55208333334; IsDate=$true}; @{Col='G'; Value=168.8999938964844; IsDate=$false}; @{Col='H'; Value=45421.38541666666; IsDate=$true}; @{Col='I'; Value=166.8500061035156; IsDate=$false}; @{Col='J'; Value="High"; IsDate=$false}; @{Col='K'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=67; Cells=@(@{Col='A'; Value="SWSOLAR.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45440.38541666666; IsDate=$true}; @{Col='E'; Value=787.8499755859375; IsDate=$false}; @{Col='F'; Value=45446.38541666666; IsDate=$true}; @{Col='G'; Value=736.9000244140625; IsDate=$false}; @{Col='H'; Value=45446.46875; IsDate=$true}; @{Col='I'; Value=736.9000244140625; IsDate=$false}; @{Col='J'; Value="High"; IsDate=$false}; @{Col='K'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=68; Cells=@(@{Col='A'; Value="SWSOLAR.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45440.38541666666; IsDate=$true}; @{Col='E'; Value=787.8499755859375; IsDate=$false}; @{Col='F'; Value=45446.42708333334; IsDate=$true}; @{Col='G'; Value=736.9000244140625; IsDate=$false}; @{Col='H'; Value=45446.46875; IsDate=$true}; @{Col='I'; Value=736.9000244140625; IsDate=$false}; @{Col='J'; Value="High"; IsDate=$false}; @{Col='K'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=69; Cells=@(@{Col='A'; Value="PCJEWELLER.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45420.38541666666; IsDate=$true}; @{Col='E'; Value=52.0; IsDate=$false}; @{Col='F'; Value=45426.42708333334; IsDate=$true}; @{Col='G'; Value=51.5; IsDate=$false}; @{Col='H'; Value=45429.38541666666; IsDate=$true}; @{Col='I'; Value=51.0; IsDate=$false}; @{Col='J'; Value="High"; IsDate=$false}; @{Col='K'; Value="06/06/2024 09:27:49"; IsDate=$false})}
)
foreach ($r in $newRows) {
    foreach ($c in $r.Cells) {
        $cell = $ws.Range($c.Col + $r.Row)
        $cell.Value = $c.Value
        if ($c.IsDate) {
            $cell.NumberFormat = $DATE_FMT
        }
    }
}

# --- two_line: append rows 17-28 ---
$ws = $wb.Worksheets("two_line")
$newRows = @(
    @{Row=17; Cells=@(@{Col='A'; Value="EICHERMOT.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45426.42708333334; IsDate=$true}; @{Col='E'; Value=4741.0; IsDate=$false}; @{Col='F'; Value=45433.38541666666; IsDate=$true}; @{Col='G'; Value=4728.60009765625; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=18; Cells=@(@{Col='A'; Value="SBILIFE.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45428.59375; IsDate=$true}; @{Col='E'; Value=1458.0; IsDate=$false}; @{Col='F'; Value=45435.51041666666; IsDate=$true}; @{Col='G'; Value=1452.949951171875; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=19; Cells=@(@{Col='A'; Value="SBILIFE.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45428.59375; IsDate=$true}; @{Col='E'; Value=1458.0; IsDate=$false}; @{Col='F'; Value=45435.59375; IsDate=$true}; @{Col='G'; Value=1452.949951171875; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=20; Cells=@(@{Col='A'; Value="NTPC.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45404.59375; IsDate=$true}; @{Col='E'; Value=342.0; IsDate=$false}; @{Col='F'; Value=45421.59375; IsDate=$true}; @{Col='G'; Value=344.8500061035156; IsDate=$false}; @{Col='H'; Value="Low"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=21; Cells=@(@{Col='A'; Value="IDFCFIRSTB.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45428.46875; IsDate=$true}; @{Col='E'; Value=76.5999984741211; IsDate=$false}; @{Col='F'; Value=45434.55208333334; IsDate=$true}; @{Col='G'; Value=76.75; IsDate=$false}; @{Col='H'; Value="Low"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=22; Cells=@(@{Col='A'; Value="AUBANK.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45427.38541666666; IsDate=$true}; @{Col='E'; Value=648.7000122070312; IsDate=$false}; @{Col='F'; Value=45442.38541666666; IsDate=$true}; @{Col='G'; Value=653.9000244140625; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=23; Cells=@(@{Col='A'; Value="ZOMATO.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45441.42708333334; IsDate=$true}; @{Col='E'; Value=185.3999938964844; IsDate=$false}; @{Col='F'; Value=45446.38541666666; IsDate=$true}; @{Col='G'; Value=184.8000030517578; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=24; Cells=@(@{Col='A'; Value="GUJGASLTD.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45419.38541666666; IsDate=$true}; @{Col='E'; Value=566.7000122070312; IsDate=$false}; @{Col='F'; Value=45434.38541666666; IsDate=$true}; @{Col='G'; Value=568.0; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=25; Cells=@(@{Col='A'; Value="SOBHA.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45415.38541666666; IsDate=$true}; @{Col='E'; Value=1967.0; IsDate=$false}; @{Col='F'; Value=45419.38541666666; IsDate=$true}; @{Col='G'; Value=1964.800048828125; IsDate=$false}; @{Col='H'; Value="High"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=26; Cells=@(@{Col='A'; Value="TV18BRDCST.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45434.38541666666; IsDate=$true}; @{Col='E'; Value=42.65000152587891; IsDate=$false}; @{Col='F'; Value=45441.38541666666; IsDate=$true}; @{Col='G'; Value=42.29999923706055; IsDate=$false}; @{Col='H'; Value="Low"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=27; Cells=@(@{Col='A'; Value="TV18BRDCST.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45428.59375; IsDate=$true}; @{Col='E'; Value=42.59999847412109; IsDate=$false}; @{Col='F'; Value=45441.38541666666; IsDate=$true}; @{Col='G'; Value=42.29999923706055; IsDate=$false}; @{Col='H'; Value="Low"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=28; Cells=@(@{Col='A'; Value="ZENTEC.NS"; IsDate=$false}; @{Col='B'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='C'; Value="hour"; IsDate=$false}; @{Col='D'; Value=45441.38541666666; IsDate=$true}; @{Col='E'; Value=917.4500122070312; IsDate=$false}; @{Col='F'; Value=45443.38541666666; IsDate=$true}; @{Col='G'; Value=916.5499877929688; IsDate=$false}; @{Col='H'; Value="Low"; IsDate=$false}; @{Col='I'; Value="06/06/2024 09:27:49"; IsDate=$false})}
)
foreach ($r in $newRows) {
    foreach ($c in $r.Cells) {
        $cell = $ws.Range($c.Col + $r.Row)
        $cell.Value = $c.Value
        if ($c.IsDate) {
            $cell.NumberFormat = $DATE_FMT
        }
    }
}

# --- ph_pl_breakout_line: append rows 232-241 ---
$ws = $wb.Worksheets("ph_pl_breakout_line")
$newRows = @(
    @{Row=232; Cells=@(@{Col='A'; Value="BHARTIARTL.NS"; IsDate=$false}; @{Col='B'; Value=45433.38541666666; IsDate=$true}; @{Col='C'; Value=1362.75; IsDate=$false}; @{Col='D'; Value=1345.599975585938; IsDate=$false}; @{Col='E'; Value=1347.25; IsDate=$false}; @{Col='F'; Value="High"; IsDate=$false}; @{Col='G'; Value=1362.75; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=1366.900024414062; IsDate=$false}; @{Col='K'; Value=1358.699951171875; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=233; Cells=@(@{Col='A'; Value="BOSCHLTD.NS"; IsDate=$false}; @{Col='B'; Value=45429.46875; IsDate=$true}; @{Col='C'; Value=30773.25; IsDate=$false}; @{Col='D'; Value=30544.69921875; IsDate=$false}; @{Col='E'; Value=30628.900390625; IsDate=$false}; @{Col='F'; Value="Low"; IsDate=$false}; @{Col='G'; Value=30544.69921875; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=30544.599609375; IsDate=$false}; @{Col='K'; Value=30610.25; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=234; Cells=@(@{Col='A'; Value="BOSCHLTD.NS"; IsDate=$false}; @{Col='B'; Value=45434.55208333334; IsDate=$true}; @{Col='C'; Value=30843.55078125; IsDate=$false}; @{Col='D'; Value=30587.55078125; IsDate=$false}; @{Col='E'; Value=30839.19921875; IsDate=$false}; @{Col='F'; Value="Low"; IsDate=$false}; @{Col='G'; Value=30587.55078125; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=30544.599609375; IsDate=$false}; @{Col='K'; Value=30610.25; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=235; Cells=@(@{Col='A'; Value="JBMA.NS"; IsDate=$false}; @{Col='B'; Value=45440.38541666666; IsDate=$true}; @{Col='C'; Value=2148.800048828125; IsDate=$false}; @{Col='D'; Value=2057.0; IsDate=$false}; @{Col='E'; Value=2069.85009765625; IsDate=$false}; @{Col='F'; Value="High"; IsDate=$false}; @{Col='G'; Value=2148.800048828125; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=2150.0; IsDate=$false}; @{Col='K'; Value=2141.699951171875; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=236; Cells=@(@{Col='A'; Value="RVNL.NS"; IsDate=$false}; @{Col='B'; Value=45443.42708333334; IsDate=$true}; @{Col='C'; Value=373.0; IsDate=$false}; @{Col='D'; Value=367.2000122070312; IsDate=$false}; @{Col='E'; Value=372.3999938964844; IsDate=$false}; @{Col='F'; Value="Low"; IsDate=$false}; @{Col='G'; Value=367.2000122070312; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=365.4500122070312; IsDate=$false}; @{Col='K'; Value=367.5499877929688; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=237; Cells=@(@{Col='A'; Value="PARAS.NS"; IsDate=$false}; @{Col='B'; Value=45443.38541666666; IsDate=$true}; @{Col='C'; Value=942.25; IsDate=$false}; @{Col='D'; Value=891.4000244140625; IsDate=$false}; @{Col='E'; Value=893.9000244140625; IsDate=$false}; @{Col='F'; Value="Low"; IsDate=$false}; @{Col='G'; Value=891.4000244140625; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=890.3499755859375; IsDate=$false}; @{Col='K'; Value=892.0999755859375; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=238; Cells=@(@{Col='A'; Value="MOTILALOFS.NS"; IsDate=$false}; @{Col='B'; Value=45439.38541666666; IsDate=$true}; @{Col='C'; Value=2304.949951171875; IsDate=$false}; @{Col='D'; Value=2243.800048828125; IsDate=$false}; @{Col='E'; Value=2264.75; IsDate=$false}; @{Col='F'; Value="High"; IsDate=$false}; @{Col='G'; Value=2304.949951171875; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=2305.0; IsDate=$false}; @{Col='K'; Value=2303.89990234375; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=239; Cells=@(@{Col='A'; Value="TATATECH.NS"; IsDate=$false}; @{Col='B'; Value=45429.59375; IsDate=$true}; @{Col='C'; Value=1047.849975585938; IsDate=$false}; @{Col='D'; Value=1044.0; IsDate=$false}; @{Col='E'; Value=1047.0; IsDate=$false}; @{Col='F'; Value="Low"; IsDate=$false}; @{Col='G'; Value=1044.0; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=1042.050048828125; IsDate=$false}; @{Col='K'; Value=1044.099975585938; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=240; Cells=@(@{Col='A'; Value="TEJASNET.NS"; IsDate=$false}; @{Col='B'; Value=45427.38541666666; IsDate=$true}; @{Col='C'; Value=1204.699951171875; IsDate=$false}; @{Col='D'; Value=1155.0; IsDate=$false}; @{Col='E'; Value=1174.449951171875; IsDate=$false}; @{Col='F'; Value="Low"; IsDate=$false}; @{Col='G'; Value=1155.0; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=1151.050048828125; IsDate=$false}; @{Col='K'; Value=1155.099975585938; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})},
    @{Row=241; Cells=@(@{Col='A'; Value="APLAPOLLO.NS"; IsDate=$false}; @{Col='B'; Value=45446.38541666666; IsDate=$true}; @{Col='C'; Value=1573.849975585938; IsDate=$false}; @{Col='D'; Value=1525.75; IsDate=$false}; @{Col='E'; Value=1533.599975585938; IsDate=$false}; @{Col='F'; Value="High"; IsDate=$false}; @{Col='G'; Value=1573.849975585938; IsDate=$false}; @{Col='H'; Value="hour"; IsDate=$false}; @{Col='I'; Value="06-06-2024 14:15:00"; IsDate=$false}; @{Col='J'; Value=1576.449951171875; IsDate=$false}; @{Col='K'; Value=1572.599975585938; IsDate=$false}; @{Col='L'; Value="06/06/2024 09:27:49"; IsDate=$false})}
)
foreach ($r in $newRows) {
    foreach ($c in $r.Cells) {
        $cell = $ws.Range($c.Col + $r.Row)
        $cell.Value = $c.Value
        if ($c.IsDate) {
            $cell.NumberFormat = $DATE_FMT
        }
    }
}
